# Updated code and protocol
# Update the "quantity" (column B) counts for several letters and
# recompute the "periodicity" (column C) share for every row as
# B<n> / SUM(B2:B34), matching the refreshed totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quantities for rows 2..34 (row 1 is the header).
$quantities = @{
    2  = 115667
    3  = 91837
    4  = 81729
    5  = 68872
    6  = 65574
    7  = 64768
    8  = 54652
    9  = 48417
    10 = 47693
    11 = 40107
    12 = 33391
    13 = 32127
    14 = 31566
    15 = 27827
    16 = 26970
    17 = 24527
    18 = 23156
    19 = 19630
    20 = 19295
    21 = 17831
    22 = 17825
    23 = 17200
    24 = 12116
    25 = 10105
    26 = 8395
    27 = 7477
    28 = 6147
    29 = 3674
    30 = 3024
    31 = 2979
    32 = 1877
    33 = 836
    34 = 307
}

# Write the refreshed quantities into column B.
foreach ($row in $quantities.Keys) {
    $ws.Cells.Item($row, 2).Value = $quantities[$row]
}

# Recompute the total and refresh the periodicity column (C) so it stays
# consistent with quantity / total(quantity).
$total = 0
foreach ($row in $quantities.Keys) {
    $total = $total + $quantities[$row]
}

foreach ($row in $quantities.Keys) {
    $ws.Cells.Item($row, 3).Value = $quantities[$row] / $total
}
